# regen save_data to use K instead of Strike#, regen std/mean, calc and write s_vals
# Update column G ("K") values for the affected rows on the active sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @{
    2  = 0
    3  = 0
    4  = 1
    5  = 0
    6  = 0
    7  = 4
    8  = 1
    9  = 0
    10 = 1
    11 = 0
    12 = 2
    13 = 1
    14 = 0
    15 = 0
    16 = 1
    17 = 1
    18 = 0
    19 = 2
    20 = 0
    22 = 1
    23 = 2
    24 = 1
    25 = 2
    26 = 1
    27 = 1
    28 = 2
    29 = 2
    30 = 0
    31 = 0
    32 = 0
    33 = 0
    34 = 1
    37 = 2
    38 = 1
    39 = 2
    40 = 1
    41 = 0
    42 = 2
    43 = 0
    44 = 1
    45 = 1
    46 = 0
    47 = 2
    49 = 1
    50 = 1
    51 = 1
    52 = 1
    54 = 1
}

foreach ($row in $values.Keys) {
    $ws.Range("G$row").Value = $values[$row]
}
